# Natmi following Dr Hou advice
# Updates the Sema5a-Plxnb3 LR-pairs sheet: the ligand/receptor specificity
# figures were recomputed (new "ECs" sending/target cluster added, all metric
# columns recalculated) -- this rewrites the full A2:T13 data block to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the 16 numeric expression/specificity metrics.
$rows = @(
    @("ECs", "Sema5a", "Plxnb3", "ECs", 2, 0.6666666666666666, 0.5290683333333334, 1.587205, 0.01267142171338989, 0.01267142171338989, 3, 1, 0.6053406666666667, 1.816022, 0.1736372177305929, 0.1736372177305929, 0.3202665776122223, 2.88239919851, 0.002200230411004042, 0.002200230411004042),
    @("ECs", "Sema5a", "Plxnb3", "FAPs", 2, 0.6666666666666666, 0.5290683333333334, 1.587205, 0.01267142171338989, 0.01267142171338989, 3, 1, 0.5337323333333334, 1.601197, 0.153096929507777, 0.153096929507777, 0.2823808760427778, 2.541427884385, 0.001939955756818166, 0.001939955756818166),
    @("ECs", "Sema5a", "Plxnb3", "M2", 2, 0.6666666666666666, 0.5290683333333334, 1.587205, 0.01267142171338989, 0.01267142171338989, 3, 1, 0.9433016666666667, 2.829905, 0.270578677263763, 0.270578677263763, 0.499071040613889, 4.491639365525, 0.003428616526260361, 0.003428616526260361),
    @("ECs", "Sema5a", "Plxnb3", "sCs", 2, 0.6666666666666666, 0.5290683333333334, 1.587205, 0.01267142171338989, 0.01267142171338989, 3, 1, 1.403863333333333, 4.21159, 0.4026871754978671, 0.4026871754978671, 0.7427396339944445, 6.68465670595, 0.005102619019307317, 0.005102619019307317),
    @("FAPs", "Sema5a", "Plxnb3", "ECs", 3, 1, 32.04971933333334, 96.149158, 0.7676050216609417, 0.7676050216609416, 3, 1, 0.6053406666666667, 1.816022, 0.1736372177305929, 0.1736372177305929, 19.40099846771956, 174.608986209476, 0.1332848002772374, 0.1332848002772374),
    @("FAPs", "Sema5a", "Plxnb3", "FAPs", 3, 1, 32.04971933333334, 96.149158, 0.7676050216609417, 0.7676050216609416, 3, 1, 0.5337323333333334, 1.601197, 0.153096929507777, 0.153096929507777, 17.10597148245845, 153.953743342126, 0.1175179718910408, 0.1175179718910408),
    @("FAPs", "Sema5a", "Plxnb3", "M2", 3, 1, 32.04971933333334, 96.149158, 0.7676050216609417, 0.7676050216609416, 3, 1, 0.9433016666666667, 2.829905, 0.270578677263763, 0.270578677263763, 30.23255366333223, 272.09298296999, 0.2076975514220397, 0.2076975514220397),
    @("FAPs", "Sema5a", "Plxnb3", "sCs", 3, 1, 32.04971933333334, 96.149158, 0.7676050216609417, 0.7676050216609416, 3, 1, 1.403863333333333, 4.21159, 0.4026871754978671, 0.4026871754978671, 44.99342581569112, 404.94083234122, 0.3091046980706237, 0.3091046980706237),
    @("sCs", "Sema5a", "Plxnb3", "ECs", 3, 1, 9.174090999999999, 27.522273, 0.2197235566256685, 0.2197235566256685, 3, 1, 0.6053406666666667, 1.816022, 0.1736372177305929, 0.1736372177305929, 5.553450362000667, 49.98105325800601, 0.03815218704235145, 0.03815218704235146),
    @("sCs", "Sema5a", "Plxnb3", "FAPs", 3, 1, 9.174090999999999, 27.522273, 0.2197235566256685, 0.2197235566256685, 3, 1, 0.5337323333333334, 1.601197, 0.153096929507777, 0.153096929507777, 4.896508995642333, 44.068580960781, 0.03363900185991801, 0.03363900185991801),
    @("sCs", "Sema5a", "Plxnb3", "M2", 3, 1, 9.174090999999999, 27.522273, 0.2197235566256685, 0.2197235566256685, 3, 1, 0.9433016666666667, 2.829905, 0.270578677263763, 0.270578677263763, 8.653935330451667, 77.88541797406499, 0.05945250931546291, 0.05945250931546291),
    @("sCs", "Sema5a", "Plxnb3", "sCs", 3, 1, 9.174090999999999, 27.522273, 0.2197235566256685, 0.2197235566256685, 3, 1, 1.403863333333333, 4.21159, 0.4026871754978671, 0.4026871754978671, 12.87916997156333, 115.91252974407, 0.08847985840793611, 0.08847985840793612)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
